$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -7
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 4
$ws.Range("F10").Value = -1
